$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 5, shifting existing rows 5-14 down to 6-15
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = "Freileitungsmast"
$ws.Cells.Item(5, 2).Value = 1251
$ws.Cells.Item(5, 3).Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/industrial"

# Copy the hyperlink cell style (the "Link" cellXf) from an existing
# hyperlinked cell onto the new C5 cell
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(5, 3).PasteSpecial(-4122)  # xlPasteFormats

# Remove the old multi-cell hyperlink that used to span C3:C14 (it no
# longer lines up correctly now that a row has been inserted in the
# middle of it) - only the C2 hyperlink remains afterwards.
$ws.Hyperlinks.Item(2).Delete()

$ws.Range("C5").Select() | Out-Null
